$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the Lapso (date range) header
$ws.Range("A1").Value = "Lapso: 2022-10-24 al 2022-11-22"

# Swap the names in rows 3 and 4, and set the "Cantidad Finalizados" counts
$ws.Range("A3").Value = "Alberto Chinsky"
$ws.Range("B3").Value = 2

$ws.Range("A4").Value = "Stefania Beatriz Marco"
$ws.Range("B4").Value = 2
